$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are plain text in the source data (some contain
# multiple "." separators, e.g. "27.563.67"), so force text formatting before
# assigning the value to stop Excel from auto-converting numeric-looking
# strings (e.g. "0.537") into floating point numbers. Reset the style back to
# Normal afterwards so no stray style/number-format is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.563.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.640.42'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.537'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.98'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.02%  '
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.10'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.634.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("E14").Value = '  -1.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.564'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.17'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.494.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0724'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.60%  '
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.55%  '
$ws.Range("E24").Value = '  -3.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.113'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.82%  '
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.55%  '
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0487'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.01%  '
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.427.09'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.45%  '
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.880'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.43%  '
$ws.Range("E39").Value = '  -2.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.879'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.10%  '
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.46%  '
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("E45").Value = '  +1.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.91'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.780.84'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.52%  '
$ws.Range("E50").Value = '  +0.75%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0989'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.34%  '
